# RCC.xlsx update
# - Fix/extend several Jira-id lists in the "Test Cases" sheet (column B)
# - Widen column B, increase row 6 height
# - Move selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Jira id (column B) text fixes -----------------------------------
$ws.Range("B4").Value  = "OPQA-1532||OPQA-1550||OPQA-3365||OPQA-3366||OPQA-1494"
$ws.Range("B6").Value  = "OPQA-3455||OPQA-3459||OPQA-3463||OPQA-1479||OPQA-1470"
$ws.Range("B12").Value = "OPQA-1448||OPQA-1451||OPQA-1454||OPQA-1464||OPQA-1465||OPQA-1462||OPQA-1457||OPQA-1459"
$ws.Range("B35").Value = "OPQA-1566||OPQA-3451||OPQA-1568"
$ws.Range("B38").Value = "OPQA-1534||OPQA-1544||OPQA-1546||OPQA-3427||OPQA-1526||OPQA-3433||OPQA-3434||OPQA-3435||OPQA-1507"
$ws.Range("B39").Value = "OPQA-1556||OPQA-3438||OPQA-1557||OPQA-1527||OPQA-3442||OPQA-3443||OPQA-3444||OPQA-1551||OPQA-1552||OPQA-1511"
$ws.Range("B40").Value = "OPQA-1529||OPQA-3445||OPQA-3446||OPQA-3447||OPQA-1564||OPQA-3450||OPQA-1565||OPQA-1513"

# --- Formatting --------------------------------------------------------
# Column B width 22.5703125 -> 50 (xml "width" = ColumnWidth + 5/6 for this sheet's font)
$ws.Columns.Item(2).ColumnWidth = 49.166666666666664

# Row 6 height 45 -> 60
$ws.Rows.Item(6).RowHeight = 60

# --- View / selection ----------------------------------------------------
$ws.Activate()
$ws.Range("B4").Select()
